$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 398, pushing the existing rows 398:430 down to 399:431
$ws.Rows("398:398").Insert()

# Populate the new row 398 with the new price-report entry
$ws.Range("A398").Value = 10
$ws.Range("B398").Value = "Vega Modelo de Temuco"
$ws.Range("C398").Value = "La Araucanía"
$ws.Range("D398").Value = 44461
$ws.Range("E398").Value = 9
$ws.Range("F398").Value = 100114001
$ws.Range("G398").Value = "Papa"
$ws.Range("H398").Value = "Rosara"
$ws.Range("I398").Value = "1a (guarda)"
$ws.Range("J398").Value = 300
$ws.Range("K398").Value = 8000
$ws.Range("L398").Value = 8000
$ws.Range("M398").Value = 8000
$ws.Range("N398").Value = "$/malla 25 kilos"
$ws.Range("O398").Value = "Provincia de Cautín"
$ws.Range("P398").Value = 320
$ws.Range("Q398").Value = 25
$ws.Range("R398").Value = "Hortaliza"
